# Labelled 75% of the data: fill in the "topic" column (D) for rows 331-434
# with the classification labels that were determined for each comment.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$values = @{
    331 = 'بطاقات'
    332 = 'حسابات'
    333 = 'مساعدة'
    334 = 'مساعدة'
    335 = 'قروض'
    336 = 'قروض'
    337 = 'قروض'
    338 = 'مساعدة'
    339 = 'مساعدة'
    340 = 'مساعدة'
    341 = 'بطاقات'
    342 = 'بطاقات'
    343 = 'بطاقات'
    344 = 'بطاقات'
    345 = 'قروض'
    346 = 'بطاقات'
    347 = 'قروض'
    348 = 'قروض'
    349 = 'مساعدة'
    350 = 'معلومة'
    351 = 'شهادات'
    352 = 'شهادات'
    353 = 'شهادات'
    354 = 'قروض'
    355 = 'مساعدة'
    356 = 'مساعدة'
    357 = 'مساعدة'
    358 = 'مساعدة'
    359 = 'مساعدة'
    360 = 'مساعدة'
    361 = 'مساعدة'
    362 = 'مساعدة'
    363 = 'قروض'
    364 = 'حسابات'
    365 = 'حسابات'
    366 = 'حسابات'
    367 = 'حسابات'
    368 = 'حسابات'
    369 = 'بطاقات'
    370 = 'بطاقات'
    371 = 'بطاقات'
    372 = 'بطاقات'
    373 = 'بطاقات'
    374 = 'مساعدة'
    375 = 'مساعدة'
    376 = 'حسابات'
    377 = 'حسابات'
    378 = 'حسابات'
    379 = 'حسابات'
    380 = 'حسابات'
    381 = 'حسابات'
    382 = 'قروض'
    383 = 'قروض'
    384 = 'قروض'
    385 = 'قروض'
    386 = 'قروض'
    387 = 'قروض'
    388 = 'قروض'
    389 = 'مساعدة'
    390 = 'مساعدة'
    391 = 'مساعدة'
    392 = 'تحويلات'
    393 = 'تحويلات'
    394 = 'معلومة'
    395 = 'قروض'
    396 = 'قروض'
    397 = 'قروض'
    398 = 'قروض'
    399 = 'قروض'
    400 = 'معلومة'
    401 = 'معلومة'
    402 = 'مساعدة'
    403 = 'مساعدة'
    404 = 'مساعدة'
    405 = 'مساعدة'
    406 = 'مساعدة'
    407 = 'مساعدة'
    408 = 'مساعدة'
    409 = 'مساعدة'
    410 = 'مساعدة'
    411 = 'قروض'
    412 = 'قروض'
    413 = 'قروض'
    414 = 'بطاقات'
    415 = 'بطاقات'
    416 = 'بطاقات'
    417 = 'بطاقات'
    418 = 'بطاقات'
    419 = 'بطاقات'
    420 = 'بطاقات'
    421 = 'قروض'
    422 = 'قروض'
    423 = 'قروض'
    424 = 'قروض'
    425 = 'قروض'
    426 = 'قروض'
    427 = 'بطاقات'
    428 = 'بطاقات'
    429 = 'بطاقات'
    430 = 'قروض'
    431 = 'قروض'
    432 = 'معلومة'
    433 = 'عقاري'
    434 = 'عقاري'
}

foreach ($r in $values.Keys) {
    $ws.Cells.Item([int]$r, 4).Value = $values[$r]
}

# Update sheet view (scroll position and selection) to match where the
# author was working when they saved the file.
$excel.ActiveWindow.ScrollRow = 425
$ws.Range("H433").Select()